$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at E:F (shifting the existing "fantasy points" column
# from E to G, carrying along its header text, style and values).
$ws.Range("E1:F1").EntireColumn.Insert()

# New header cells (reuse E1's/headers formatting automatically via insert)
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Fill in height/weight for every data row (2-17)
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 5).Value = 6.5
    $ws.Cells.Item($row, 6).Value = 265
}
